$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 302-303, pushing existing rows (302..363) down to (304..365).
$ws.Range("A302:A303").EntireRow.Insert()

# Row 302: new "Primera" quality record for date 2022-03-21 (serial 44641)
$ws.Cells.Item(302, 1).Value = 8
$ws.Cells.Item(302, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(302, 3).Value = "Coquimbo"
$ws.Cells.Item(302, 4).Value = 44641
$ws.Cells.Item(302, 5).Value = 4
$ws.Cells.Item(302, 6).Value = 100112009
$ws.Cells.Item(302, 7).Value = "Acelga"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 2500
$ws.Cells.Item(302, 11).Value = 500
$ws.Cells.Item(302, 12).Value = 600
$ws.Cells.Item(302, 13).Value = 550
$ws.Cells.Item(302, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(302, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(302, 16).Value = 275
$ws.Cells.Item(302, 17).Value = 2
$ws.Cells.Item(302, 18).Value = "Hortaliza"

# Row 303: new "Segunda" quality record for the same date
$ws.Cells.Item(303, 1).Value = 8
$ws.Cells.Item(303, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(303, 3).Value = "Coquimbo"
$ws.Cells.Item(303, 4).Value = 44641
$ws.Cells.Item(303, 5).Value = 4
$ws.Cells.Item(303, 6).Value = 100112009
$ws.Cells.Item(303, 7).Value = "Acelga"
$ws.Cells.Item(303, 8).Value = "Sin especificar"
$ws.Cells.Item(303, 9).Value = "Segunda"
$ws.Cells.Item(303, 10).Value = 1300
$ws.Cells.Item(303, 11).Value = 400
$ws.Cells.Item(303, 12).Value = 450
$ws.Cells.Item(303, 13).Value = 425
$ws.Cells.Item(303, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(303, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(303, 16).Value = 212
$ws.Cells.Item(303, 17).Value = 2
$ws.Cells.Item(303, 18).Value = "Hortaliza"
